$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.371.70'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '2.353.26'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.33'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.91'
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('D9').Value = '2.370.71'
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('E11').Value = '  +5.69%  '
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.344'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.22'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').Value = '2.775.89'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').Value = '57.379.52'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '2.365.76'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '329.39'
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.24'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.33'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('E25').Value = '  +4.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.23'
$ws.Range('E27').Value = '  +7.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.32'
$ws.Range('E28').Value = '  +11.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.99'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = '0.0₃0746'
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.29'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.63'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('E35').Value = '  +2.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.925'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.05'
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('E39').Value = '  +4.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.52'
$ws.Range('E40').Value = '  +2.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '150.95'
$ws.Range('E41').Value = '  +8.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.385'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  +2.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.31'
$ws.Range('E44').Value = '  +4.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '282.24'
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0939'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0508'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.565'
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0220'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.11'
$ws.Range('E50').Value = '  +6.05%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.63'
$ws.Range('E51').Value = '  +4.21%  '
